$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.567.62"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.018.92"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.40"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.601"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.87"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.39"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0748"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.315.15"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.23"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.14"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.764"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.11"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.008.94"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.453.89"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0797"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.39"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.10"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.54%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.12"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.62"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.06%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +3.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "18.86"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.36"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0600"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.23"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.32"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.69"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.40%  "
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.455.29"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.20"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +40.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0203"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.03"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  -6.92%  "
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.86"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.36%  "
